$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column E (imputed/restored values) for several rows ---
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("E17").Value = -7.3
$ws.Range("E18").Value = -8.5
$ws.Range("E19").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("E23").Value = -7

# --- Remove two whole rows (RM 232 and SC 92) ---
# Row 26 is "RM 232"; deleting it shifts everything below up by one.
$ws.Rows("26").Delete()
# After the above deletion, the row that used to be "SC 92" (originally row 28)
# is now row 27.
$ws.Rows("27").Delete()

# --- Apply the remaining cell-level edits on the now-shifted bottom block ---
$ws.Range("B27").Value = -20.4
$ws.Range("E27").ClearContents()

$ws.Range("B28").ClearContents()

$ws.Range("B29").ClearContents()

$ws.Range("B30").Value = -19.7

$ws.Range("B32").ClearContents()
